$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F7").Value = 2778
$ws1.Range("F8").Value = 1668
$ws1.Range("F9").Value = 1759
$ws1.Range("F13").Value = 881
$ws1.Range("F14").Value = 159
$ws1.Range("F20").Value = 6398
$ws1.Range("F22").Value = 1398
$ws1.Range("F24").Value = 175
$ws1.Range("F25").Value = 156
$ws1.Range("F26").Value = 299
$ws1.Range("F27").Value = 256
$ws1.Range("F28").Value = 61
$ws1.Range("F34").Value = 462
$ws1.Range("F35").Value = 1332
$ws1.Range("F38").Value = 215
$ws1.Range("F40").Value = 141
$ws1.Range("F41").Value = 180

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 11

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F7").Value = 11
$ws4.Range("F10").Value = 2778
$ws4.Range("F11").Value = 1668
$ws4.Range("F12").Value = 1759
$ws4.Range("F17").Value = 881
$ws4.Range("F18").Value = 159
$ws4.Range("F23").Value = 6398
$ws4.Range("F25").Value = 1398
$ws4.Range("F28").Value = 175
$ws4.Range("F29").Value = 156
$ws4.Range("F30").Value = 299
$ws4.Range("F31").Value = 256
$ws4.Range("F32").Value = 61
$ws4.Range("F38").Value = 462
$ws4.Range("F39").Value = 1332
$ws4.Range("F42").Value = 215
$ws4.Range("F44").Value = 141
$ws4.Range("F45").Value = 180
